$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $rows = $used.Rows.Count
    $cols = $used.Columns.Count

    for ($r = 1; $r -le $rows; $r++) {
        for ($c = 1; $c -le $cols; $c++) {
            $cell = $used.Cells.Item($r, $c)
            $v = $cell.Value2
            if ($v -ne $null -and $v -is [string]) {
                # Strip footnote-style bracket annotations like "[1]", "[2]", "[5, 6]"
                $new = [System.Text.RegularExpressions.Regex]::Replace($v, "\[[0-9, ]+\]", "")
                # Collapse embedded line breaks into a single space
                $new = $new.Replace("`r`n", "`n")
                $new = $new.Replace("`n", " ")
                if ($new -ne $v) {
                    $cell.Value = $new
                }
            }
        }
    }
}
